$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update lastUSD (column B) and Dif (column C) values per the fixed
# binary tree calculation.

$ws.Range("B2").Value = 2686.97731952
$ws.Range("C2").Value = 0

$ws.Range("B3").Value = 2709.1707084
$ws.Range("C3").Value = 0.008259611541479028

$ws.Range("B4").Value = 2968.538408
$ws.Range("C4").Value = 0.1047872962806764

$ws.Range("B5").Value = 3057.67
$ws.Range("C5").Value = 0.1379589912378645

$ws.Range("B6").Value = 3058.56
$ws.Range("C6").Value = 0.1382902184475376

$ws.Range("B7").Value = 3070.7
$ws.Range("C7").Value = 0.1428083064536427
